$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width (new custom width ~21 chars) ---
$ws.Range("B1").ColumnWidth = 20.1

# --- Rebuild hyperlinks: K2's old hyperlink goes away; K3-K8 keep links but
#     K5-K8 now point at https:// targets instead of http:// ---
$ws.Hyperlinks.Delete()

# K2 is plain text now (no hyperlink), pointing at a brand new stream URL
$ws.Range("K2").Value = "https://streaming.megaguardiao.com.br/click190-cloud/3e026423e27946aba02e.stream/chunklist_w1960725029.m3u8"

# K3, K4 keep the same targets/text as before
$ws.Range("K3").Value = "https://nginx.megaguardiao.com.br/live/hls/m-lojagil.m3u8"
$ws.Hyperlinks.Add($ws.Range("K3"), "https://nginx.megaguardiao.com.br/live/hls/m-lojagil.m3u8")
$ws.Range("K4").Value = "https://nginx.megaguardiao.com.br/live/hls/m-ozeaisgaragem.m3u8"
$ws.Hyperlinks.Add($ws.Range("K4"), "https://nginx.megaguardiao.com.br/live/hls/m-ozeaisgaragem.m3u8")

# K5-K8 flip from http:// to https:// (update displayed text + link target)
$ws.Range("K5").Value = "https://nginx.megaguardiao.com.br/live/hls/m-mbtosa1.m3u8"
$ws.Hyperlinks.Add($ws.Range("K5"), "https://nginx.megaguardiao.com.br/live/hls/m-mbtosa1.m3u8")
$ws.Range("K6").Value = "https://nginx.megaguardiao.com.br/live/hls/m-mbtosa2.m3u8"
$ws.Hyperlinks.Add($ws.Range("K6"), "https://nginx.megaguardiao.com.br/live/hls/m-mbtosa2.m3u8")
$ws.Range("K7").Value = "https://nginx.megaguardiao.com.br/live/hls/mbtosa3.m3u8"
$ws.Hyperlinks.Add($ws.Range("K7"), "https://nginx.megaguardiao.com.br/live/hls/mbtosa3.m3u8")
$ws.Range("K8").Value = "https://nginx.megaguardiao.com.br/live/hls/lpr.m3u8"
$ws.Hyperlinks.Add($ws.Range("K8"), "https://nginx.megaguardiao.com.br/live/hls/lpr.m3u8")

# Restore the original (non-duplicated) hyperlink cell style that Add() nudges
$ws.Range("K3:K8").Style = "Hiperlink"

# --- New row 16: a lone formatted (underlined), empty cell at K16 ---
$ws.Range("K16").Font.Underline = 1

# --- Final selection cursor used when the workbook was last saved ---
$ws.Range("J10").Select()
